$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '69.301.57'
$ws.Range("E2").Value = '  -0.62%  '

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '3.848.48'
$ws.Range("E3").Value = '  +4.26%  '

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$dCell = $ws.Range("D4")
$dCell.Value = "'1.00"
$dCell.Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$dCell = $ws.Range("D5")
$dCell.Value = "'603.73"
$dCell.Style = "Normal"
$ws.Range("E5").Value = '  -1.78%  '

$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$dCell = $ws.Range("D6")
$dCell.Value = "'173.12"
$dCell.Style = "Normal"
$ws.Range("E6").Value = '  -2.73%  '

$ws.Range("B7").Value = 'LidoStakedEther'
$ws.Range("C7").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D7").Value = '3.846.20'
$ws.Range("E7").Value = '  +4.21%  '

$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$dCell = $ws.Range("D8")
$dCell.Value = "'1.00"
$dCell.Style = "Normal"
$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("B9").Value = 'XRP'
$ws.Range("C9").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$dCell = $ws.Range("D9")
$dCell.Value = "'0.527"
$dCell.Style = "Normal"
$ws.Range("E9").Value = '  -0.43%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$dCell = $ws.Range("D10")
$dCell.Value = "'0.164"
$dCell.Style = "Normal"
$ws.Range("E10").Value = '  -0.03%  '

$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$dCell = $ws.Range("D11")
$dCell.Value = "'6.38"
$dCell.Style = "Normal"
$ws.Range("E11").Value = '  +2.33%  '

$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$dCell = $ws.Range("D12")
$dCell.Value = "'0.478"
$dCell.Style = "Normal"
$ws.Range("E12").Value = '  -0.11%  '

$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$dCell = $ws.Range("D13")
$dCell.Value = "'39.28"
$dCell.Style = "Normal"
$ws.Range("E13").Value = '  -0.98%  '

$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$dCell = $ws.Range("D14")
$dCell.Value = "'0.0000251"
$dCell.Style = "Normal"
$ws.Range("E14").Value = '  -0.45%  '

$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '4.485.05'
$ws.Range("E15").Value = '  +4.19%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.841.06'
$ws.Range("E16").Value = '  +3.98%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '69.452.27'
$ws.Range("E17").Value = '  -0.39%  '

$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$dCell = $ws.Range("D18")
$dCell.Value = "'7.44"
$dCell.Style = "Normal"
$ws.Range("E18").Value = '  -0.95%  '

$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$dCell = $ws.Range("D19")
$dCell.Value = "'0.117"
$dCell.Style = "Normal"
$ws.Range("E19").Value = '  -3.54%  '

$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$dCell = $ws.Range("D20")
$dCell.Value = "'16.39"
$dCell.Style = "Normal"
$ws.Range("E20").Value = '  +0.49%  '

$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$dCell = $ws.Range("D21")
$dCell.Value = "'501.29"
$dCell.Style = "Normal"
$ws.Range("E21").Value = '  +0.26%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$dCell = $ws.Range("D22")
$dCell.Value = "'9.58"
$dCell.Style = "Normal"
$ws.Range("E22").Value = '  +4.93%  '

$ws.Range("B23").Value = 'Polygon'
$ws.Range("C23").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$dCell = $ws.Range("D23")
$dCell.Value = "'0.746"
$dCell.Style = "Normal"
$ws.Range("E23").Value = '  +4.93%  '

$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$dCell = $ws.Range("D24")
$dCell.Value = "'87.52"
$dCell.Style = "Normal"
$ws.Range("E24").Value = '  +1.68%  '

$ws.Range("B25").Value = 'Fetch.AI'
$ws.Range("C25").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$dCell = $ws.Range("D25")
$dCell.Value = "'2.41"
$dCell.Style = "Normal"
$ws.Range("E25").Value = '  -4.32%  '

$ws.Range("B26").Value = 'PEPE'
$ws.Range("C26").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$dCell = $ws.Range("D26")
$dCell.Value = "'0.0000138"
$dCell.Style = "Normal"
$ws.Range("E26").Value = '  +8.40%  '

$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$dCell = $ws.Range("D27")
$dCell.Value = "'12.54"
$dCell.Style = "Normal"
$ws.Range("E27").Value = '  -2.84%  '

$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$dCell = $ws.Range("D28")
$dCell.Value = "'10.28"
$dCell.Style = "Normal"
$ws.Range("E28").Value = '  -9.47%  '

$ws.Range("B29").Value = 'Dai'
$ws.Range("C29").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$dCell = $ws.Range("D29")
$dCell.Value = "'1.00"
$dCell.Style = "Normal"
$ws.Range("E29").Value = '  +0.14%  '

$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$dCell = $ws.Range("D30")
$dCell.Value = "'2.51"
$dCell.Style = "Normal"
$ws.Range("E30").Value = '  +3.48%  '

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$dCell = $ws.Range("D31")
$dCell.Value = "'2.97"
$dCell.Style = "Normal"
$ws.Range("E31").Value = '  +3.27%  '

$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$dCell = $ws.Range("D32")
$dCell.Value = "'33.10"
$dCell.Style = "Normal"
$ws.Range("E32").Value = '  +10.33%  '

$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$dCell = $ws.Range("D33")
$dCell.Value = "'7.85"
$dCell.Style = "Normal"
$ws.Range("E33").Value = '  -1.07%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$dCell = $ws.Range("D34")
$dCell.Value = "'0.113"
$dCell.Style = "Normal"
$ws.Range("E34").Value = '  -0.10%  '

$ws.Range("B35").Value = 'FirstDigitalUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$dCell = $ws.Range("D35")
$dCell.Value = "'0.999"
$dCell.Style = "Normal"
$ws.Range("E35").Value = '  -0.07%  '

$ws.Range("B36").Value = 'Mantle'
$ws.Range("C36").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$dCell = $ws.Range("D36")
$dCell.Value = "'1.04"
$dCell.Style = "Normal"
$ws.Range("E36").Value = '  -0.90%  '

$ws.Range("B37").Value = 'Filecoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$dCell = $ws.Range("D37")
$dCell.Value = "'6.03"
$dCell.Style = "Normal"
$ws.Range("E37").Value = '  +0.20%  '

$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$dCell = $ws.Range("D38")
$dCell.Value = "'0.139"
$dCell.Style = "Normal"
$ws.Range("E38").Value = '  +2.11%  '

$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$dCell = $ws.Range("D39")
$dCell.Value = "'463.88"
$dCell.Style = "Normal"
$ws.Range("E39").Value = '  +8.68%  '

$ws.Range("B40").Value = 'TheGraph'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$dCell = $ws.Range("D40")
$dCell.Value = "'0.331"
$dCell.Style = "Normal"
$ws.Range("E40").Value = '  -1.47%  '

$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$dCell = $ws.Range("D41")
$dCell.Value = "'2.05"
$dCell.Style = "Normal"
$ws.Range("E41").Value = '  -0.06%  '

$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$dCell = $ws.Range("D42")
$dCell.Value = "'49.54"
$dCell.Style = "Normal"
$ws.Range("E42").Value = '  -0.80%  '

$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$dCell = $ws.Range("D43")
$dCell.Value = "'2.87"
$dCell.Style = "Normal"
$ws.Range("E43").Value = '  -0.99%  '

$ws.Range("B44").Value = 'Cosmos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$dCell = $ws.Range("D44")
$dCell.Value = "'8.51"
$dCell.Style = "Normal"
$ws.Range("E44").Value = '  -0.14%  '

$ws.Range("B45").Value = 'Arweave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$dCell = $ws.Range("D45")
$dCell.Value = "'42.33"
$dCell.Style = "Normal"
$ws.Range("E45").Value = '  -4.44%  '

$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '2.894.31'
$ws.Range("E46").Value = '  -1.43%  '

$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$dCell = $ws.Range("D47")
$dCell.Value = "'0.0359"
$dCell.Style = "Normal"
$ws.Range("E47").Value = '  +0.55%  '

$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$dCell = $ws.Range("D48")
$dCell.Value = "'27.51"
$dCell.Style = "Normal"
$ws.Range("E48").Value = '  +1.10%  '

$ws.Range("B49").Value = 'USDe'
$ws.Range("C49").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$dCell = $ws.Range("D49")
$dCell.Value = "'1.00"
$dCell.Style = "Normal"
$ws.Range("E49").Value = '  -0.01%  '

$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$dCell = $ws.Range("D50")
$dCell.Value = "'139.14"
$dCell.Style = "Normal"
$ws.Range("E50").Value = '  +2.14%  '

$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$dCell = $ws.Range("D51")
$dCell.Value = "'2.38"
$dCell.Style = "Normal"
$ws.Range("E51").Value = '  -1.50%  '
